$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 1.516666666666667
$ws.Range("C7").Value = 964.8
$ws.Range("D7").Value = 967
$ws.Range("B14").Value = 2.556603773584906
$ws.Range("C14").Value = 651.5877192982456
$ws.Range("B19").Value = 2.193548387096774
$ws.Range("C19").Value = 924.6764705882352
$ws.Range("E19").Value = 34
$ws.Range("B34").Value = 2.256410256410256
$ws.Range("C34").Value = 824.5348837209302
$ws.Range("D34").Value = 1256
$ws.Range("E34").Value = 44
$ws.Range("B47").Value = 2.051020408163265
$ws.Range("C47").Value = 1080.707070707071
$ws.Range("D47").Value = 1045.625
$ws.Range("E47").Value = 99
$ws.Range("B56").Value = 2.043478260869565
$ws.Range("C56").Value = 1123.075
$ws.Range("D56").Value = 1233.290322580645
$ws.Range("B58").Value = 1.929203539823009
$ws.Range("C58").Value = 1461.016666666667
$ws.Range("D58").Value = 1077.533333333333
$ws.Range("B61").Value = 2.008620689655173
$ws.Range("C61").Value = 1242.873949579832
$ws.Range("D61").Value = 1096.290322580645
$ws.Range("B63").Value = 1.957983193277311
$ws.Range("C63").Value = 1067.075
$ws.Range("D63").Value = 1345.212765957447
$ws.Range("C64").Value = 853.9416666666667
$ws.Range("D64").Value = 756.8518518518518
$ws.Range("B67").Value = 2.255102040816327
$ws.Range("C67").Value = 1045.15
$ws.Range("D67").Value = 927.2048192771084
$ws.Range("B68").Value = 2.92929292929293
$ws.Range("C68").Value = 1715.308333333333
$ws.Range("D68").Value = 1298.787878787879
$ws.Range("B72").Value = 1.567567567567568
$ws.Range("C72").Value = 1117.756302521008
$ws.Range("D72").Value = 990.0566037735849
$ws.Range("B73").Value = 1.522935779816514
$ws.Range("C73").Value = 864.2
$ws.Range("D73").Value = 786.8139534883721
$ws.Range("B75").Value = 2.258064516129032
$ws.Range("C75").Value = 1100.290322580645
$ws.Range("E75").Value = 31
$ws.Range("B76").Value = 1.831858407079646
$ws.Range("C76").Value = 1010.15
$ws.Range("D76").Value = 1183.381818181818
$ws.Range("B77").Value = 1.445454545454546
$ws.Range("C77").Value = 1467.716666666667
$ws.Range("D77").Value = 839.5434782608696
$ws.Range("B79").Value = 2.355769230769231
$ws.Range("C79").Value = 1463.766666666667
$ws.Range("D79").Value = 961.9111111111112
$ws.Range("B81").Value = 1.909090909090909
$ws.Range("C81").Value = 1267.268907563025
$ws.Range("B83").Value = 1.756521739130435
$ws.Range("C83").Value = 1072.15
$ws.Range("D83").Value = 938.2962962962963
$ws.Range("B85").Value = 1.481818181818182
$ws.Range("C85").Value = 1380.866666666667
$ws.Range("D85").Value = 887.8039215686274
$ws.Range("B87").Value = 2.071428571428572
$ws.Range("C87").Value = 904.4074074074074
$ws.Range("E87").Value = 81
$ws.Range("B89").Value = 1.588235294117647
$ws.Range("C89").Value = 812.6875
$ws.Range("D89").Value = 1047.857142857143
$ws.Range("E89").Value = 64
$ws.Range("C90").Value = 944.3666666666667
$ws.Range("D90").Value = 831.1375
$ws.Range("B93").Value = 1.425
$ws.Range("C93").Value = 1235.375
$ws.Range("D93").Value = 857.4347826086956
$ws.Range("B95").Value = 3.105263157894737
$ws.Range("C95").Value = 1774.808333333333
